$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Rename header row.
#    Columns A:J carried the "_old" suffix -> rename to "_FV2404".
#    Columns L:U carried the "_new" suffix -> rename to "_FV2410".
#    Column K ("diff") is left untouched.
# ---------------------------------------------------------------------------
$oldSuffixHeaders = @(
    "Segmentname",
    "Segmentgruppe",
    "Segment",
    "Datenelement",
    "Segment ID",
    "Code",
    "Qualifier",
    "Beschreibung",
    "Bedingungsausdruck",
    "Bedingung"
)

for ($i = 0; $i -lt $oldSuffixHeaders.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = "$($oldSuffixHeaders[$i])_FV2404"
}

for ($i = 0; $i -lt $oldSuffixHeaders.Length; $i++) {
    $ws.Cells.Item(1, 11 + $i + 1).Value = "$($oldSuffixHeaders[$i])_FV2410"
}

# ---------------------------------------------------------------------------
# 2. Convert the used range into a real Excel Table (ListObject) named
#    "Table1", covering A1:U77 with headers.
#
#    The header row (A1:U1) already carries bold/fill/border/center/wrap
#    formatting applied directly to the cells. If that formatting is left in
#    place while the table is created, Excel captures it into a table-level
#    "header row" differential format (headerRowDxfId) on top of the quick
#    style. To keep the header's direct cell formatting as the only source of
#    its look (matching the original authoring), temporarily stash the
#    formatting, strip it before creating the table, then restore it
#    afterwards.
# ---------------------------------------------------------------------------
$headerRange = $ws.Range("A1:U1")
$scratchRange = $ws.Range("A1000:U1000")

$headerRange.Copy()
$scratchRange.PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$headerRange.ClearFormats()

$listRange = $ws.Range("A1:U77")
$table = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $listRange, $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$table.Name = "Table1"

$scratchRange.Copy()
$headerRange.PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Drop the scratch row entirely so it doesn't linger as an empty formatted
# row and inflate the sheet's used range / dimension.
$ws.Rows.Item(1000).Delete()

# ---------------------------------------------------------------------------
# 3. Freeze the header row (split after row 1).
# ---------------------------------------------------------------------------
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
